# Weekly refresh: the values in columns D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg) get reshuffled across the data rows
# (rows 2-23, row 10 is untouched). This mirrors a re-sort of the
# underlying daily records into a new weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values currently sitting in the source
# row before this script runs should end up in the destination row)
$mapping = @{
    2  = 23
    3  = 7
    4  = 2
    5  = 22
    6  = 21
    7  = 3
    8  = 13
    9  = 5
    11 = 8
    12 = 19
    13 = 9
    14 = 6
    15 = 17
    16 = 18
    17 = 11
    18 = 4
    19 = 12
    20 = 15
    21 = 14
    22 = 16
    23 = 20
}

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# First, snapshot the "before" values of every relevant cell so the
# subsequent writes (which happen in-place) never read already-overwritten
# data.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range($c + $row).Value2
    }
    $snapshot[$row] = $rowData
}

# Now write the snapshot from each source row into its destination row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range($c + $destRow).Value = $rowData[$c]
    }
}
